$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.419.43"
$ws.Range("E2").Value = "  +1.00%  "

$ws.Range("D3").Value = "1.852.17"
$ws.Range("E3").Value = "  +1.29%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4745"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2757"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06346"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.40%  "

$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.824.42"
$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07470"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.978"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6235"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.05%  "

$ws.Range("D16").Value = "30.376.97"
$ws.Range("E16").Value = "  +1.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "245.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007350"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9996"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.929"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.903"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "164.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.044"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.874"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1024"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.347"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.039"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.831"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04826"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.128"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.6976"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.701"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01893"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.682"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8767"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.988"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "106.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4065"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.497"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.175"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1200"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.531"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05501"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.349"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3685"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.80%  "
